# Edit script: insert 3 new data rows into the "Feria Lagunitas de Puerto Montt - Uva" sheet
# at position 240 (pushing existing rows 240-318 down to 243-321), then populate the
# 3 new rows with their values (dimension grows from A1:T318 to A1:T321).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 240, shifting rows 240:318 down to 243:321.
$ws.Rows("240:242").Insert()

# Common (constant) values shared by every data row in this block.
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100109
$producto  = "Uva"
$categoriaId = 100109001
$categoria = "Uva"
$calidad   = "Primera"

# New row 240: Red Globe
$r = 240
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44985
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Red Globe"
$ws.Cells.Item($r, 12).Value = $calidad
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 14000
$ws.Cells.Item($r, 15).Value = 15000
$ws.Cells.Item($r, 16).Value = 14500
$ws.Cells.Item($r, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 806
$ws.Cells.Item($r, 20).Value = 18

# New row 241: Superior Seedless
$r = 241
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44985
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Superior Seedless"
$ws.Cells.Item($r, 12).Value = $calidad
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 14000
$ws.Cells.Item($r, 15).Value = 15000
$ws.Cells.Item($r, 16).Value = 14500
$ws.Cells.Item($r, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 806
$ws.Cells.Item($r, 20).Value = 18

# New row 242: Thompson seedless
$r = 242
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44985
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Thompson seedless"
$ws.Cells.Item($r, 12).Value = $calidad
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 15000
$ws.Cells.Item($r, 15).Value = 15500
$ws.Cells.Item($r, 16).Value = 15250
$ws.Cells.Item($r, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 847
$ws.Cells.Item($r, 20).Value = 18
